$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename product
$ws.Range("A2").Value = "Salchipapas"

# Row 3: update cost/PVP/ganancia
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1

# Row 4: rename product and update PVP/ganancia
$ws.Range("A4").Value = "gorros"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 2

# Row 5: delete entirely (Chochos con tostado)
$ws.Rows.Item(5).Delete()
